# Refresh the cryptocurrency price/volume table (coinranking.com feed).
# Each entry below is one cell whose text changed in this update; cells
# whose new value looks like a plain number are first formatted as Text
# so Excel keeps the exact original digits instead of parsing them into a
# floating point number (matching the inline-string cells already in the
# sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "35.492.68"; ForceText = $false },
    @{ Cell = "E2"; Value = "  +2.91%  "; ForceText = $false },
    @{ Cell = "D3"; Value = "1.852.20"; ForceText = $false },
    @{ Cell = "E3"; Value = "  +2.40%  "; ForceText = $false },
    @{ Cell = "E4"; Value = "  +0.54%  "; ForceText = $false },
    @{ Cell = "D5"; Value = "229.69"; ForceText = $true },
    @{ Cell = "E5"; Value = "  +1.92%  "; ForceText = $false },
    @{ Cell = "D6"; Value = "0.608"; ForceText = $true },
    @{ Cell = "E6"; Value = "  +3.51%  "; ForceText = $false },
    @{ Cell = "E7"; Value = "  +0.44%  "; ForceText = $false },
    @{ Cell = "D8"; Value = "41.72"; ForceText = $true },
    @{ Cell = "E8"; Value = "  +9.40%  "; ForceText = $false },
    @{ Cell = "E9"; Value = "  +6.43%  "; ForceText = $false },
    @{ Cell = "D10"; Value = "0.0691"; ForceText = $true },
    @{ Cell = "E10"; Value = "  +2.68%  "; ForceText = $false },
    @{ Cell = "E11"; Value = "  +3.95%  "; ForceText = $false },
    @{ Cell = "D12"; Value = "2.121.20"; ForceText = $false },
    @{ Cell = "E12"; Value = "  +2.48%  "; ForceText = $false },
    @{ Cell = "D13"; Value = "11.43"; ForceText = $true },
    @{ Cell = "E13"; Value = "  +2.79%  "; ForceText = $false },
    @{ Cell = "D14"; Value = "1.852.96"; ForceText = $false },
    @{ Cell = "E14"; Value = "  +2.17%  "; ForceText = $false },
    @{ Cell = "D15"; Value = "0.670"; ForceText = $true },
    @{ Cell = "E15"; Value = "  +6.61%  "; ForceText = $false },
    @{ Cell = "D16"; Value = "4.68"; ForceText = $true },
    @{ Cell = "E16"; Value = "  +6.06%  "; ForceText = $false },
    @{ Cell = "D17"; Value = "35.519.16"; ForceText = $false },
    @{ Cell = "E17"; Value = "  +3.09%  "; ForceText = $false },
    @{ Cell = "D18"; Value = "69.85"; ForceText = $true },
    @{ Cell = "E18"; Value = "  +2.70%  "; ForceText = $false },
    @{ Cell = "D19"; Value = "247.16"; ForceText = $true },
    @{ Cell = "E19"; Value = "  +2.04%  "; ForceText = $false },
    @{ Cell = "D20"; Value = "0.0₃0798"; ForceText = $false },
    @{ Cell = "E20"; Value = "  +3.74%  "; ForceText = $false },
    @{ Cell = "D21"; Value = "12.08"; ForceText = $true },
    @{ Cell = "E21"; Value = "  +8.68%  "; ForceText = $false },
    @{ Cell = "D22"; Value = "4.62"; ForceText = $true },
    @{ Cell = "E22"; Value = "  +12.54%  "; ForceText = $false },
    @{ Cell = "D23"; Value = "1.01"; ForceText = $true },
    @{ Cell = "E23"; Value = "  +0.45%  "; ForceText = $false },
    @{ Cell = "E24"; Value = "  -0.53%  "; ForceText = $false },
    @{ Cell = "D25"; Value = "168.95"; ForceText = $true },
    @{ Cell = "E25"; Value = "  -0.56%  "; ForceText = $false },
    @{ Cell = "D26"; Value = "7.90"; ForceText = $true },
    @{ Cell = "E26"; Value = "  +2.41%  "; ForceText = $false },
    @{ Cell = "D27"; Value = "17.74"; ForceText = $true },
    @{ Cell = "E27"; Value = "  +1.08%  "; ForceText = $false },
    @{ Cell = "E28"; Value = "  +1.47%  "; ForceText = $false },
    @{ Cell = "E29"; Value = "  +13.01%  "; ForceText = $false },
    @{ Cell = "E30"; Value = "  +0.55%  "; ForceText = $false },
    @{ Cell = "D31"; Value = "3.311.85"; ForceText = $false },
    @{ Cell = "E31"; Value = "  +36.31%  "; ForceText = $false },
    @{ Cell = "E32"; Value = "  +5.49%  "; ForceText = $false },
    @{ Cell = "D33"; Value = "3.91"; ForceText = $true },
    @{ Cell = "E33"; Value = "  +3.89%  "; ForceText = $false },
    @{ Cell = "D34"; Value = "4.05"; ForceText = $true },
    @{ Cell = "E34"; Value = "  +5.66%  "; ForceText = $false },
    @{ Cell = "D35"; Value = "1.86"; ForceText = $true },
    @{ Cell = "E35"; Value = "  +2.76%  "; ForceText = $false },
    @{ Cell = "D36"; Value = "98.67"; ForceText = $true },
    @{ Cell = "E36"; Value = "  +20.84%  "; ForceText = $false },
    @{ Cell = "D37"; Value = "0.695"; ForceText = $true },
    @{ Cell = "E37"; Value = "  +8.65%  "; ForceText = $false },
    @{ Cell = "D38"; Value = "1.362.12"; ForceText = $false },
    @{ Cell = "E38"; Value = "  +1.40%  "; ForceText = $false },
    @{ Cell = "B39"; Value = "TrustWalletToken"; ForceText = $false },
    @{ Cell = "C39"; Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"; ForceText = $false },
    @{ Cell = "D39"; Value = "1.08"; ForceText = $true },
    @{ Cell = "E39"; Value = "  +2.56%  "; ForceText = $false },
    @{ Cell = "B40"; Value = "RenderToken"; ForceText = $false },
    @{ Cell = "C40"; Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; ForceText = $false },
    @{ Cell = "D40"; Value = "2.45"; ForceText = $true },
    @{ Cell = "E40"; Value = "  +5.33%  "; ForceText = $false },
    @{ Cell = "E41"; Value = "  +3.45%  "; ForceText = $false },
    @{ Cell = "D42"; Value = "1.00"; ForceText = $true },
    @{ Cell = "E42"; Value = "  +6.09%  "; ForceText = $false },
    @{ Cell = "E43"; Value = "  +4.26%  "; ForceText = $false },
    @{ Cell = "B44"; Value = "HuobiToken"; ForceText = $false },
    @{ Cell = "C44"; Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; ForceText = $false },
    @{ Cell = "D44"; Value = "2.48"; ForceText = $true },
    @{ Cell = "E44"; Value = "  +1.17%  "; ForceText = $false },
    @{ Cell = "B45"; Value = "InjectiveProtocol"; ForceText = $false },
    @{ Cell = "C45"; Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; ForceText = $false },
    @{ Cell = "D45"; Value = "14.63"; ForceText = $true },
    @{ Cell = "E45"; Value = "  +7.11%  "; ForceText = $false },
    @{ Cell = "E46"; Value = "  +0.84%  "; ForceText = $false },
    @{ Cell = "D47"; Value = "0.0520"; ForceText = $true },
    @{ Cell = "E47"; Value = "  +1.70%  "; ForceText = $false },
    @{ Cell = "D48"; Value = "6.20"; ForceText = $true },
    @{ Cell = "E48"; Value = "  +8.13%  "; ForceText = $false },
    @{ Cell = "D49"; Value = "2.019.51"; ForceText = $false },
    @{ Cell = "E49"; Value = "  +2.49%  "; ForceText = $false },
    @{ Cell = "E50"; Value = "  +0.38%  "; ForceText = $false },
    @{ Cell = "D51"; Value = "103.66"; ForceText = $true },
    @{ Cell = "E51"; Value = "  +1.43%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Pre-format as Text so a numeric-looking string (e.g. "41.72")
        # is stored verbatim instead of becoming a Double.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        # Drop back to the default style now that the text is committed,
        # so the cell does not carry a lingering custom number format.
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
